$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "43.815.68"
$ws.Range("E2").Value = "  +0.37%  "

# Row 3
$ws.Range("D3").Value = "2.315.63"
$ws.Range("E3").Value = "  +1.23%  "

# Row 4
$ws.Range("E4").Value = "  +0.22%  "

# Row 5
$ws.Range("D5").Value = "'112.40"
$ws.Range("E5").Value = "  +17.59%  "

# Row 6
$ws.Range("D6").Value = "'270.63"
$ws.Range("E6").Value = "  +1.32%  "

# Row 7
$ws.Range("D7").Value = "'0.621"
$ws.Range("E7").Value = "  -0.51%  "

# Row 8
$ws.Range("E8").Value = "  +0.18%  "

# Row 9
$ws.Range("D9").Value = "'0.624"
$ws.Range("E9").Value = "  +2.44%  "

# Row 10
$ws.Range("E10").Value = "  +7.11%  "

# Row 11
$ws.Range("D11").Value = "'0.0943"
$ws.Range("E11").Value = "  +0.69%  "

# Row 12
$ws.Range("D12").Value = "'8.90"
$ws.Range("E12").Value = "  +15.25%  "

# Row 13
$ws.Range("E13").Value = "  +2.37%  "

# Row 14
$ws.Range("D14").Value = "'15.77"
$ws.Range("E14").Value = "  +3.94%  "

# Row 15
$ws.Range("D15").Value = "2.663.45"
$ws.Range("E15").Value = "  +1.74%  "

# Row 16
$ws.Range("D16").Value = "'0.862"
$ws.Range("E16").Value = "  +1.47%  "

# Row 17
$ws.Range("D17").Value = "2.320.50"

# Row 18
$ws.Range("D18").Value = "43.889.35"
$ws.Range("E18").Value = "  +0.70%  "

# Row 19
$ws.Range("E19").Value = "  +2.29%  "

# Row 20
$ws.Range("D20").Value = "'6.69"
$ws.Range("E20").Value = "  +8.59%  "

# Row 21
$ws.Range("D21").Value = "'72.58"
$ws.Range("E21").Value = "  +0.40%  "

# Row 22
$ws.Range("D22").Value = "'2.50"
$ws.Range("E22").Value = "  +6.79%  "

# Row 23
$ws.Range("D23").Value = "'234.25"
$ws.Range("E23").Value = "  -0.37%  "

# Row 24
$ws.Range("D24").Value = "'9.52"
$ws.Range("E24").Value = "  +6.73%  "

# Row 25
$ws.Range("E25").Value = "  +15.75%  "

# Row 26
$ws.Range("E26").Value = "  +0.01%  "

# Row 27
$ws.Range("D27").Value = "'11.56"
$ws.Range("E27").Value = "  +2.18%  "

# Row 28
$ws.Range("D28").Value = "'42.57"
$ws.Range("E28").Value = "  +12.03%  "

# Row 29
$ws.Range("E29").Value = "  +0.03%  "

# Row 31
$ws.Range("D31").Value = "'177.96"
$ws.Range("E31").Value = "  +1.17%  "

# Row 32
$ws.Range("E32").Value = "  -0.10%  "

# Row 33
$ws.Range("D33").Value = "'0.0935"
$ws.Range("E33").Value = "  +5.71%  "

# Row 34
$ws.Range("D34").Value = "'5.62"
$ws.Range("E34").Value = "  +4.62%  "

# Row 35
$ws.Range("E35").Value = "  +1.32%  "

# Row 36
$ws.Range("D36").Value = "'4.79"
$ws.Range("E36").Value = "  +8.35%  "

# Row 37
$ws.Range("D37").Value = "'0.112"
$ws.Range("E37").Value = "  +3.67%  "

# Row 38
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.0360"
$ws.Range("E38").Value = "  +1.73%  "

# Row 39
$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D39").Value = "'3.93"
$ws.Range("E39").Value = "  +20.46%  "

# Row 40
$ws.Range("D40").Value = "'0.247"
$ws.Range("E40").Value = "  +4.36%  "

# Row 41
$ws.Range("D41").Value = "'2.39"
$ws.Range("E41").Value = "  +0.74%  "

# Row 42
$ws.Range("D42").Value = "'69.60"
$ws.Range("E42").Value = "  +11.77%  "

# Row 43
$ws.Range("B43").Value = "Celestia"
$ws.Range("C43").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D43").Value = "'12.93"
$ws.Range("E43").Value = "  +8.79%  "

# Row 44
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  +0.08%  "

# Row 45
$ws.Range("D45").Value = "'1.40"
$ws.Range("E45").Value = "  +4.21%  "

# Row 46
$ws.Range("E46").Value = "  +10.48%  "

# Row 47
$ws.Range("D47").Value = "'8.83"
$ws.Range("E47").Value = "  +0.11%  "

# Row 48
$ws.Range("E48").Value = "  -1.02%  "

# Row 49
$ws.Range("D49").Value = "'100.10"
$ws.Range("E49").Value = "  +1.96%  "

# Row 50
$ws.Range("B50").Value = "WOONetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D50").Value = "'0.465"
$ws.Range("E50").Value = "  +11.13%  "

# Row 51
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").Value = "'1.22"
$ws.Range("E51").Value = "  +2.65%  "

Write-Host "done"